# LMS-2523 Update BaSynthec Validation
# The "Strain" value on the openbis-metadata sheet (cell B3) is updated
# from "MGP9" to "JJS-MGP9".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")
$ws.Range("B3").Value = "JJS-MGP9"
